$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A, B and E (COM ColumnWidth uses character units that are
# offset from the stored XML width by 5/6 of a character in this engine,
# so compensate to land on the exact target widths of 18 / 15 / 17).
$ws.Columns.Item(1).ColumnWidth = 17.166666666666668
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 16.166666666666668

# Shipping details to append as four identical new rows (5-8)
$address = "123 Maple Street"
$city = "Beverly Hills"
$firstName = "Emily"
$lastName = "Johnson"
$phone = "+1-310-555-0199"
$state = "California"
$zip = "90210"

for ($r = 5; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $address
    $ws.Cells.Item($r, 2).Value = $city
    $ws.Cells.Item($r, 3).Value = $firstName
    $ws.Cells.Item($r, 4).Value = $lastName
    $ws.Cells.Item($r, 5).Value = $phone
    $ws.Cells.Item($r, 6).Value = $state

    # Force the zip code to be stored as text (not a number) while keeping
    # the cell's style at the default (no explicit style), matching the
    # original unstyled data rows.
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $zip
    $ws.Cells.Item($r, 7).Style = "Normal"
}
